# !Maturita.xlsx - "Farma zvířat, Bylo nás pět"
#
# The reading checklist on sheet "List1" is updated: the "Doporučené"
# (recommended / "D") marker is moved off of a few rows and a book that
# has now actually been read gets its checkbox filled in with 1.
#
#   C44 (Čapek Karel - Povídky z jedné kapsy, Povídky z druhé kapsy): D -> (blank)
#   C45 (Čapek Karel - Válka s mloky):                                (blank) -> D
#   C51 (Hrabal Bohumil - Ostře sledované vlaky):                     (blank) -> D
#   C53 (Jirotka Zdeněk - Saturnin):                                  D -> (blank)
#   C61 (Smoliak/Svěrák - Dobytí severního pólu):                     D -> 1 (read!)
#   C67 (Wolker Jiří - Host do domu):                                 D -> (blank)
#   C68 (Wolker Jiří - Těžká hodina):                                 D -> (blank)
#
# All of the SUM/COUNTIF totals further down the sheet (C69, D69, C71,
# D71, C72) are plain formulas and recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List1")
$ws.Activate()

$ws.Range("C44").ClearContents()
$ws.Range("C45").Value = "D"
$ws.Range("C51").Value = "D"
$ws.Range("C53").ClearContents()
$ws.Range("C61").Value = 1
$ws.Range("C67").ClearContents()
$ws.Range("C68").ClearContents()

# Restore the view: keep row 1 frozen, scroll the body down a bit and
# leave the selection on B15.
$excel.ActiveWindow.ScrollRow = 26
$ws.Range("B15").Select()
